$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update D2 value from "demo" to "smartmed" (A2 keeps its value "CA-F77PUNMY")
$ws.Range("D2").Value = "smartmed"

# Move selection from A2 to D2
$ws.Range("D2").Select()
